# Update metric values in rows 2-26 (columns B:Q) with new values.
# Every data row (2 through 26) shares the same set of new values per column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B" = 0.3082271070678201
    "C" = -2.241727458159763
    "D" = 0.5184404671956622
    "E" = -2.424969855977629
    "F" = -0.0205656919417887
    "G" = 0.4106662722775318
    "H" = 1.924429454498304
    "I" = 0.1077442308400815
    "J" = 0.3371630501027281
    "K" = 0.2224536404714048
    "L" = 0.2690475247287556
    "M" = 0.6408324837877148
    "N" = -1.07531867879654
    "O" = 0.6681140493286761
    "P" = 33.77994876784079
    "Q" = 53.281961965732
}

for ($row = 2; $row -le 26; $row++) {
    foreach ($col in $newValues.Keys) {
        $ws.Range("$col$row").Value = $newValues[$col]
    }
}
